$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: Student ID and Log Time changed
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "231994"
$ws.Range("A2").Style = "Normal"

$ws.Range("D2").Value = "10:53:21"

# Insert a new row 3 with a new scan entry
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "231995"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").Value = "Anatomy"
$ws.Range("C3").Value = "14/08/2025"
$ws.Range("D3").Value = "10:53:21"
$ws.Range("E3").Value = "Scan"
$ws.Range("F3").Value = "admin@admin.com"
